$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column L ("Hire Date Str") is being removed: it duplicated the Hire Date
# values in column E (with a couple of values corrupted into the strings
# "physics"/"maths"). We drop the header and all the data in that column,
# but leave the now-empty cells in place (same as the row/col dimensions).
$ws.Range("L1").ClearContents()
$ws.Range("L2:L8").ClearContents()
$ws.Range("L10:L14").ClearContents()

# A couple of the old cells had a text quote-prefix on them; drop all prior
# formatting on the cleared cells before giving them their new look below.
$ws.Range("L2:L8").ClearFormats()
$ws.Range("L10:L14").ClearFormats()

# The cleared L cells pick up the same "real date" look as column E below.
$ws.Range("L2:L8").Font.Bold = $true
$ws.Range("L2:L8").NumberFormat = "mm-dd-yy"
$ws.Range("L10:L14").Font.Bold = $true
$ws.Range("L10:L14").NumberFormat = "mm-dd-yy"

# --- Column E ("Hire Date") switches from a plain integer serial-number
# display to an actual bold date format.
$ws.Range("E2:E8").Font.Bold = $true
$ws.Range("E2:E8").NumberFormat = "mm-dd-yy"
$ws.Range("E10:E14").Font.Bold = $true
$ws.Range("E10:E14").NumberFormat = "mm-dd-yy"

# The blank separator row (row 9) keeps its shaded fill in both E and L, but
# also picks up the new bold date format even though it carries no value.
$ws.Range("L9").Font.Bold = $true
$ws.Range("L9").NumberFormat = "mm-dd-yy"

# Restore the selection back to the top-left of the data, matching the
# pre-"Hire Date Str" workbook state.
$ws.Range("D13").Select()
